# Update the control scenario simulation results: new min load values (col B)
# and new min load times (col C), with two of the time strings now reused
# (duplicated) rather than unique.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newB = @(1.01134, 1.01456, 1.01123, 1.01825, 1.00783, 1.00856, 1.01328, 1.0012, 1.00574, 1.00806)
$newC = @("18:25:00", "18:28:00", "19:28:00", "18:26:00", "18:21:00", "18:33:00", "18:34:00", "18:48:00", "18:28:00", "18:16:00")

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newB[$i]
    $ws.Cells.Item($row, 3).Value = $newC[$i]
}
